$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 9 (shifts existing rows 9-28 down to 10-29,
# carrying the formatting of the row below, as Excel normally does).
$ws.Rows(9).Insert()

# Populate the new "Administrar" menu entry (child of "Restaurante"),
# matching the layout of sibling rows.
$ws.Range("A9").Value = "Administrar"
$ws.Range("B9").Value = "#"
$ws.Range("C9").Value = "Restaurante"
$ws.Range("D9").Value = 1
$ws.Range("E9").Value = "developer_board"
$ws.Range("F9").Value = "ADMIN_GERENTE"

# Re-parent the items that used to hang directly off "Restaurante" so they
# now belong to the new "Administrar" submenu.
$ws.Range("C10").Value = "Administrar"
$ws.Range("C11").Value = "Administrar"
$ws.Range("C12").Value = "Administrar"
$ws.Range("C13").Value = "Administrar"
$ws.Range("C14").Value = "Administrar"

# Underline the icon for the first re-parented row (matches the workbook as
# committed).
$ws.Range("E10").Font.Underline = $true

# Leave the selection where the author left it when saving.
$ws.Range("F9").Select()
